# Update the NIEM mapping worksheet: retarget the "Person" rows from the
# evaluation-search-results (evalsres) schema to the health-information
# search-results (phisres) schema, and rename the second business class
# from "Evaluation" to "Behavioral Health" with its new XPath.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Person / Unique and Timeless ID -> PersonPersistentIdentification (phisres)
$ws.Range("C2").Value = "/phisres-doc:PersonHealthInformationSearchResults/nc30:Person/phisres-ext:PersonPersistentIdentification"

# Row 3: Person / Temporary ID -> PersonTemporaryIdentification/IdentificationID (phisres)
$ws.Range("C3").Value = "/phisres-doc:PersonHealthInformationSearchResults/nc30:Person/phisres-ext:PersonTemporaryIdentification/nc30:IdentificationID"

# Row 4: "Evaluation" business class renamed to "Behavioral Health"
$ws.Range("A4").Value = "Behavioral Health"
$ws.Range("C4").Value = "/phisres-doc:PersonHealthInformationSearchResults/phisres-ext:BehavioralHealthInformation/jxdm51:Evaluation/jxdm51:EvaluationDiagnosisDescriptionText"

# Row 5 (SMI Indicator) text is unchanged.

# Row heights for the (now taller, wrapped) Behavioral Health rows.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30

# Column A needs to be a bit wider to fit "Behavioral Health".
$ws.Columns.Item(1).ColumnWidth = 14.83

# Selection moves to C3.
$ws.Range("C3").Select() | Out-Null
